$d = $word.ActiveDocument

# "Versi" + "on" (two separate runs) both read "Version" -- merge them into
# a single run by replacing the full word (this leaves the surrounding
# spellStart/spellEnd proofErr markers untouched).
$d.Content.Find.Execute("Version", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Version", 2) | Out-Null

# " 2" (its own run) becomes " 1." -- stop the match right before the
# bookmark so bookmarkStart/bookmarkEnd are preserved in place.
$d.Content.Find.Execute(" 2", $true, $false, $false, $false, $false,
                         $true, 1, $false, " 1.", 2) | Out-Null

# The trailing "." run (after the bookmark) is now redundant -- remove just
# that trailing character, leaving the bookmark as the last thing in the
# paragraph before the paragraph mark.
$p = $d.Paragraphs.Item(1)
$r = $d.Range($p.Range.End - 2, $p.Range.End - 1)
$r.Delete()
